$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 957.8182
$ws.Range("I6").Value = 1098.5714
$ws.Range("J6").Value = 711.5
$ws.Range("K6").Value = 3295.7142
$ws.Range("L6").Value = 2134.5
$ws.Range("M6").Value = -3183.7142
$ws.Range("N6").Value = -2358.5
$ws.Range("H33").Value = 316.7037
$ws.Range("I33").Value = 311.625
$ws.Range("K33").Value = 311.625
$ws.Range("M33").Value = -82.625
$ws.Range("H40").Value = 55558160.0
$ws.Range("J40").Value = 125002984.0
$ws.Range("L40").Value = 125002984.0
$ws.Range("N40").Value = -125003334.0
$ws.Range("H70").Value = 1525577.8
$ws.Range("I70").Value = 4065707.2
$ws.Range("K70").Value = 12197121.6
$ws.Range("M70").Value = -12196851.6
$ws.Range("H73").Value = 1525577.8
$ws.Range("I73").Value = 4065707.2
$ws.Range("K73").Value = 12197121.6
$ws.Range("M73").Value = -12196185.6
$ws.Range("H98").Value = 2720.0
$ws.Range("I98").Value = 2238.0
$ws.Range("K98").Value = 2238.0
$ws.Range("M98").Value = -740.0
$ws.Range("H122").Value = 2720.0
$ws.Range("I122").Value = 2238.0
$ws.Range("K122").Value = 6714.0
$ws.Range("M122").Value = -4264.0
$ws.Range("H125").Value = 2159.0
$ws.Range("I125").Value = 2120.6667
$ws.Range("J125").Value = 2216.5
$ws.Range("K125").Value = 19086.0003
$ws.Range("L125").Value = 19948.5
$ws.Range("M125").Value = -16626.0003
$ws.Range("N125").Value = -24868.5
$ws.Range("H129").Value = 4587.0
$ws.Range("I129").Value = 796.375
$ws.Range("J129").Value = 9641.167
$ws.Range("K129").Value = 2389.125
$ws.Range("L129").Value = 28923.501
$ws.Range("M129").Value = 2610.875
$ws.Range("N129").Value = -38923.501
$ws.Range("H131").Value = 7577077.0
$ws.Range("I131").Value = 1979.5
$ws.Range("J131").Value = 22727272.0
$ws.Range("K131").Value = 5938.5
$ws.Range("L131").Value = 68181816.0
$ws.Range("M131").Value = -898.5
$ws.Range("N131").Value = -68191896.0
$ws.Range("H137").Value = 905.875
$ws.Range("I137").Value = 905.875
$ws.Range("K137").Value = 2717.625
$ws.Range("M137").Value = -167.625
$ws.Range("H138").Value = 12069.625
$ws.Range("J138").Value = 13975.0
$ws.Range("L138").Value = 41925.0
$ws.Range("N138").Value = -52205.0
$ws.Range("H141").Value = 15157121.0
$ws.Range("I141").Value = 17860878.0
$ws.Range("K141").Value = 53582634.0
$ws.Range("M141").Value = -53577454.0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2664.1128
$ws.Range("I32").Value = 2799.9385
$ws.Range("K32").Value = 2799.9385
$ws.Range("M32").Value = -2512.9385
$ws.Range("H76").Value = 67996.336
$ws.Range("J76").Value = 67996.336
$ws.Range("L76").Value = 67996.336
$ws.Range("N76").Value = -68672.336
$ws.Range("H79").Value = 67996.336
$ws.Range("J79").Value = 67996.336
$ws.Range("L79").Value = 67996.336
$ws.Range("N79").Value = -70336.336
$ws.Range("H122").Value = 4541.75
$ws.Range("I122").Value = 3638.2856
$ws.Range("K122").Value = 10914.8568
$ws.Range("M122").Value = -8464.856800000001
$ws.Range("H132").Value = 4003485.2
$ws.Range("I132").Value = 3090.647
$ws.Range("K132").Value = 9271.940999999999
$ws.Range("M132").Value = -6741.940999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2371.5557
$ws.Range("J80").Value = 2877.1
$ws.Range("L80").Value = 2877.1
$ws.Range("N80").Value = -4873.1
$ws.Range("H83").Value = 2371.5557
$ws.Range("J83").Value = 2877.1
$ws.Range("L83").Value = 14385.5
$ws.Range("N83").Value = -24369.5
$ws.Range("H86").Value = 4372.6875
$ws.Range("I86").Value = 1611.5
$ws.Range("J86").Value = 6029.4
$ws.Range("K86").Value = 1611.5
$ws.Range("L86").Value = 6029.4
$ws.Range("M86").Value = -488.5
$ws.Range("N86").Value = -8275.4
$ws.Range("H89").Value = 4372.6875
$ws.Range("I89").Value = 1611.5
$ws.Range("J89").Value = 6029.4
$ws.Range("K89").Value = 8057.5
$ws.Range("L89").Value = 30147.0
$ws.Range("M89").Value = -2441.5
$ws.Range("N89").Value = -41379.0
$ws.Range("H103").Value = 49530.832
$ws.Range("J103").Value = 49530.832
$ws.Range("L103").Value = 49530.832
$ws.Range("N103").Value = -51874.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1152.2142
$ws.Range("I107").Value = 340.9375
$ws.Range("J107").Value = 2233.9167
$ws.Range("K107").Value = 340.9375
$ws.Range("L107").Value = 2233.9167
$ws.Range("M107").Value = 1579.0625
$ws.Range("N107").Value = -6073.9167
$ws.Range("H122").Value = 2954.24
$ws.Range("I122").Value = 3005.3845
$ws.Range("J122").Value = 2898.8333
$ws.Range("K122").Value = 9016.1535
$ws.Range("L122").Value = 8696.499899999999
$ws.Range("M122").Value = -6566.1535
$ws.Range("N122").Value = -13596.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 715.3333
$ws.Range("I98").Value = 734.5
$ws.Range("K98").Value = 2203.5
$ws.Range("M98").Value = -705.5
$ws.Range("H126").Value = 100015784.0
$ws.Range("I126").Value = 125011400.0
$ws.Range("K126").Value = 375034200.0
$ws.Range("M126").Value = -375029260.0
$ws.Range("H131").Value = 3138.9583
$ws.Range("I131").Value = 1436.2142
$ws.Range("J131").Value = 5522.8
$ws.Range("K131").Value = 4308.642599999999
$ws.Range("L131").Value = 16568.4
$ws.Range("M131").Value = 731.3574000000008
$ws.Range("N131").Value = -26648.4
$ws.Range("H132").Value = 1002.5
$ws.Range("J132").Value = 1002.5
$ws.Range("L132").Value = 9022.5
$ws.Range("N132").Value = -14082.5
$ws.Range("H140").Value = 2792.5
$ws.Range("I140").Value = 1195.5
$ws.Range("J140").Value = 7982.75
$ws.Range("K140").Value = 3586.5
$ws.Range("L140").Value = 23948.25
$ws.Range("M140").Value = 1593.5
$ws.Range("N140").Value = -34308.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7209.737
$ws.Range("I70").Value = 6557.1113
$ws.Range("K70").Value = 6557.1113
$ws.Range("M70").Value = -6287.1113
$ws.Range("H73").Value = 7209.737
$ws.Range("I73").Value = 6557.1113
$ws.Range("K73").Value = 6557.1113
$ws.Range("M73").Value = -5621.1113
$ws.Range("H80").Value = 2684.3333
$ws.Range("I80").Value = 1537.3334
$ws.Range("J80").Value = 4978.3335
$ws.Range("K80").Value = 1537.3334
$ws.Range("L80").Value = 4978.3335
$ws.Range("M80").Value = -539.3334
$ws.Range("N80").Value = -6974.3335
$ws.Range("H83").Value = 2684.3333
$ws.Range("I83").Value = 1537.3334
$ws.Range("J83").Value = 4978.3335
$ws.Range("K83").Value = 7686.666999999999
$ws.Range("L83").Value = 24891.6675
$ws.Range("M83").Value = -2694.666999999999
$ws.Range("N83").Value = -34875.6675
$ws.Range("H102").Value = 4121.0
$ws.Range("I102").Value = 4121.0
$ws.Range("K102").Value = 4121.0
$ws.Range("M102").Value = -2499.0
$ws.Range("H122").Value = 2451.5454
$ws.Range("J122").Value = 998.0
$ws.Range("L122").Value = 2994.0
$ws.Range("N122").Value = -7894.0
$ws.Range("H132").Value = 3451281.2
$ws.Range("I132").Value = 2891.4
$ws.Range("K132").Value = 8674.2
$ws.Range("M132").Value = -6144.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8335686.0
$ws.Range("I68").Value = 13890541.0
$ws.Range("K68").Value = 13890541.0
$ws.Range("M68").Value = -13889792.0
$ws.Range("H71").Value = 8335686.0
$ws.Range("I71").Value = 13890541.0
$ws.Range("K71").Value = 69452705.0
$ws.Range("M71").Value = -69448961.0
$ws.Range("H82").Value = 4070.6785
$ws.Range("I82").Value = 3613.182
$ws.Range("J82").Value = 4366.706
$ws.Range("K82").Value = 3613.182
$ws.Range("L82").Value = 4366.706
$ws.Range("M82").Value = -3252.182
$ws.Range("N82").Value = -5088.706
$ws.Range("H85").Value = 4070.6785
$ws.Range("I85").Value = 3613.182
$ws.Range("J85").Value = 4366.706
$ws.Range("K85").Value = 3613.182
$ws.Range("L85").Value = 4366.706
$ws.Range("M85").Value = -2365.182
$ws.Range("N85").Value = -6862.706
$ws.Range("H122").Value = 3535.8
$ws.Range("I122").Value = 2881.0
$ws.Range("J122").Value = 4699.8887
$ws.Range("K122").Value = 8643.0
$ws.Range("L122").Value = 14099.6661
$ws.Range("M122").Value = -6193.0
$ws.Range("N122").Value = -18999.6661
$ws.Range("H132").Value = 2032.0286
$ws.Range("I132").Value = 1941.75
$ws.Range("K132").Value = 5825.25
$ws.Range("M132").Value = -3295.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 25067.334
$ws.Range("J101").Value = 25067.334
$ws.Range("L101").Value = 25067.334
$ws.Range("N101").Value = -31557.334
$ws.Range("H107").Value = 3173.795
$ws.Range("I107").Value = 1561.9524
$ws.Range("J107").Value = 5054.278
$ws.Range("K107").Value = 4685.857199999999
$ws.Range("L107").Value = 15162.834
$ws.Range("M107").Value = -2765.857199999999
$ws.Range("N107").Value = -19002.834
$ws.Range("H122").Value = 3199.75
$ws.Range("I122").Value = 3333.0
$ws.Range("J122").Value = 2800.0
$ws.Range("K122").Value = 9999.0
$ws.Range("L122").Value = 8400.0
$ws.Range("M122").Value = -7549.0
$ws.Range("N122").Value = -13300.0
$ws.Range("H126").Value = 3319.84
$ws.Range("I126").Value = 2794.8823
$ws.Range("K126").Value = 8384.6469
$ws.Range("M126").Value = -5914.6469
$ws.Range("H137").Value = 73666.0
$ws.Range("J137").Value = 73666.0
$ws.Range("L137").Value = 73666.0
$ws.Range("N137").Value = -83866.0
